# Add daily power records: extend the comforter-cda table from row 79
# (A1:F79) through row 90 (A1:F90) with new daily readings, matching the
# existing Date / Start Time / End Time / Duration / Second Duration /
# Absolute Value layout and formulas used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 79 (existing row): fill in the previously-empty Start/End Time ---
$ws.Range("B79").Value = 0
$ws.Range("C79").Value = 0

# --- Row 80 ---
$ws.Range("A80").Value = 43404
$ws.Range("B80").Value = 0
$ws.Range("C80").Value = 0
$ws.Range("D80").Formula = "=(C80-B80)* 1440"
$ws.Range("E80").Formula = "=IF(C80>B80, (C80-B80)*1440, (B80-C80)*1440)"
$ws.Range("F80").Formula = "=ABS((C80-B80)*1440)"

# --- Row 81 ---
$ws.Range("A81").Value = 43405
$ws.Range("B81").Value = 0
$ws.Range("C81").Value = 0
$ws.Range("D81").Formula = "=(C81-B81)* 1440"
$ws.Range("E81").Formula = "=IF(C81>B81, (C81-B81)*1440, (B81-C81)*1440)"
$ws.Range("F81").Formula = "=ABS((C81-B81)*1440)"

# --- Row 82 ---
$ws.Range("A82").Value = 43406
$ws.Range("B82").Value = 0
$ws.Range("C82").Value = 0
$ws.Range("D82").Formula = "=(C82-B82)* 1440"
$ws.Range("E82").Formula = "=IF(C82>B82, (C82-B82)*1440, (B82-C82)*1440)"
$ws.Range("F82").Formula = "=ABS((C82-B82)*1440)"

# --- Row 83 ---
$ws.Range("A83").Value = 43407
$ws.Range("B83").Value = 0
$ws.Range("C83").Value = 0
$ws.Range("D83").Formula = "=(C83-B83)* 1440"
$ws.Range("E83").Formula = "=IF(C83>B83, (C83-B83)*1440, (B83-C83)*1440)"
$ws.Range("F83").Formula = "=ABS((C83-B83)*1440)"

# --- Row 84 ---
$ws.Range("A84").Value = 43408
$ws.Range("B84").Value = 0
$ws.Range("C84").Value = 0
$ws.Range("D84").Formula = "=(C84-B84)* 1440"
$ws.Range("E84").Formula = "=IF(C84>B84, (C84-B84)*1440, (B84-C84)*1440)"
$ws.Range("F84").Formula = "=ABS((C84-B84)*1440)"

# --- Row 85 ---
$ws.Range("A85").Value = 43409
$ws.Range("B85").Value = 0
$ws.Range("C85").Value = 0
$ws.Range("D85").Formula = "=(C85-B85)* 1440"
$ws.Range("E85").Formula = "=IF(C85>B85, (C85-B85)*1440, (B85-C85)*1440)"
$ws.Range("F85").Formula = "=ABS((C85-B85)*1440)"

# --- Row 86 (only a Start Time, no End Time recorded) ---
$ws.Range("A86").Value = 43410
$ws.Range("B86").Value = 0.77500000000000002
$ws.Range("D86").Formula = "=(C86-B86)* 1440"
$ws.Range("E86").Formula = "=IF(C86>B86, (C86-B86)*1440, (B86-C86)*1440)"
$ws.Range("F86").Formula = "=ABS((C86-B86)*1440)"

# --- Row 87 (no Start/End Time recorded) ---
$ws.Range("A87").Value = 43411
$ws.Range("D87").Formula = "=(C87-B87)* 1440"
$ws.Range("E87").Formula = "=IF(C87>B87, (C87-B87)*1440, (B87-C87)*1440)"
$ws.Range("F87").Formula = "=ABS((C87-B87)*1440)"

# --- Row 88 ---
$ws.Range("A88").Value = 43412
$ws.Range("D88").Formula = "=(C88-B88)* 1440"
$ws.Range("E88").Formula = "=IF(C88>B88, (C88-B88)*1440, (B88-C88)*1440)"
$ws.Range("F88").Formula = "=ABS((C88-B88)*1440)"

# --- Row 89 ---
$ws.Range("A89").Value = 43413
$ws.Range("D89").Formula = "=(C89-B89)* 1440"
$ws.Range("E89").Formula = "=IF(C89>B89, (C89-B89)*1440, (B89-C89)*1440)"
$ws.Range("F89").Formula = "=ABS((C89-B89)*1440)"

# --- Row 90 ---
$ws.Range("A90").Value = 43414
$ws.Range("D90").Formula = "=(C90-B90)* 1440"
$ws.Range("E90").Formula = "=IF(C90>B90, (C90-B90)*1440, (B90-C90)*1440)"
$ws.Range("F90").Formula = "=ABS((C90-B90)*1440)"

# --- Grow the table / autofilter range to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F90"))

# --- Update the view: scroll so the new rows are visible and select C91 ---
$ws.Range("C91").Select()
$excel.ActiveWindow.ScrollRow = 69
